$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 22:05"

# --- Simple in-place stat refreshes (country order unchanged) ---
# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 1637955
$ws.Cells.Item(4, 3).Value = 17053
$ws.Cells.Item(4, 4).Value = 395127
$ws.Cells.Item(4, 5).Value = 1145531
$ws.Cells.Item(4, 7).Value = 943
$ws.Cells.Item(4, 8).Value = 97297

# Alemania (row 11)
$ws.Cells.Item(11, 2).Value = 179713
$ws.Cells.Item(11, 3).Value = 692
$ws.Cells.Item(11, 7).Value = 43
$ws.Cells.Item(11, 8).Value = 8352

# Camerun (row 68)
$ws.Cells.Item(68, 2).Value = 4400
$ws.Cells.Item(68, 3).Value = 112
$ws.Cells.Item(68, 4).Value = 1822
$ws.Cells.Item(68, 5).Value = 2419
$ws.Cells.Item(68, 7).Value = 3
$ws.Cells.Item(68, 8).Value = 159

# Maldivas (row 101) - only D/E changed
$ws.Cells.Item(101, 4).Value = 109
$ws.Cells.Item(101, 5).Value = 1161

# --- Reordered countries: "Sudan del Sur" now overtakes Jamaica/Nepal/Tanzania ---
$ws.Cells.Item(129, 1).Value = "Sudan del Sur"
$ws.Cells.Item(129, 2).Value = 563
$ws.Cells.Item(129, 3).Value = 82
$ws.Cells.Item(129, 4).Value = 6
$ws.Cells.Item(129, 5).Value = 551
$ws.Cells.Item(129, 7).Value = 2
$ws.Cells.Item(129, 8).Value = 6

$ws.Cells.Item(130, 1).Value = "Jamaica"
$ws.Cells.Item(130, 2).Value = 534
$ws.Cells.Item(130, 3).Value = 5
$ws.Cells.Item(130, 4).Value = 181
$ws.Cells.Item(130, 5).Value = 344
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 9

$ws.Cells.Item(131, 1).Value = "Nepal"
$ws.Cells.Item(131, 2).Value = 516
$ws.Cells.Item(131, 3).Value = 59
$ws.Cells.Item(131, 4).Value = 70
$ws.Cells.Item(131, 5).Value = 443
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 3

$ws.Cells.Item(132, 1).Value = "Tanzania"
$ws.Cells.Item(132, 2).Value = 509
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 183
$ws.Cells.Item(132, 5).Value = 305
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 21

# --- "Togo" now overtakes "Cabo Verde" ---
$ws.Cells.Item(140, 1).Value = "Togo"
$ws.Cells.Item(140, 2).Value = 363
$ws.Cells.Item(140, 3).Value = 9
$ws.Cells.Item(140, 4).Value = 121
$ws.Cells.Item(140, 5).Value = 230
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 12

$ws.Cells.Item(141, 1).Value = "Cabo Verde"
$ws.Cells.Item(141, 2).Value = 362
$ws.Cells.Item(141, 3).Value = 6
$ws.Cells.Item(141, 4).Value = 95
$ws.Cells.Item(141, 5).Value = 264
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 3

# --- "Gambia" now overtakes "Antigua y Barbuda" ---
$ws.Cells.Item(189, 1).Value = "Gambia"
$ws.Cells.Item(189, 2).Value = 25
$ws.Cells.Item(189, 3).Value = 1
$ws.Cells.Item(189, 4).Value = 13
$ws.Cells.Item(189, 5).Value = 11
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 1

$ws.Cells.Item(190, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(190, 2).Value = 25
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 19
$ws.Cells.Item(190, 5).Value = 3
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 3
